# Update the division problems in the table to the new set of values.
# Each old value below is unique within the document, so a simple
# Find/Replace (wdReplaceAll) for each pair is safe. The only subtlety is
# ordering: "396÷4=" is both an existing value (to become "550÷4=") and a
# brand new value produced by another replacement ("702÷6=" -> "396÷4=").
# To avoid the second replacement accidentally altering text created by a
# later one, the pair "396÷4=" -> "550÷4=" is performed before
# "702÷6=" -> "396÷4=". All other pairs are independent and are applied in
# their natural (document) order.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "743÷3=" "960÷2="
Replace-Text "245÷7=" "743÷9="
Replace-Text "984÷5=" "784÷7="
Replace-Text "665÷9=" "890÷8="
Replace-Text "437÷3=" "312÷7="
Replace-Text "832÷5=" "657÷4="
Replace-Text "767÷3=" "631÷4="
Replace-Text "519÷4=" "880÷5="
Replace-Text "295÷4=" "763÷8="
Replace-Text "474÷4=" "755÷4="
Replace-Text "718÷4=" "219÷8="
Replace-Text "369÷6=" "245÷5="
Replace-Text "944÷8=" "569÷8="
Replace-Text "104÷9=" "882÷5="
Replace-Text "955÷4=" "574÷3="
Replace-Text "396÷4=" "550÷4="
Replace-Text "702÷6=" "396÷4="
Replace-Text "465÷4=" "223÷4="
Replace-Text "319÷4=" "553÷7="
Replace-Text "604÷3=" "980÷9="
Replace-Text "302÷6=" "331÷8="
Replace-Text "653÷5=" "293÷3="
Replace-Text "207÷8=" "685÷9="
Replace-Text "637÷7=" "139÷5="
Replace-Text "502÷9=" "590÷8="

Write-Host "All replacements applied."
